# ---------------------------------------------------------------------------
# Adds two new worksheets ("Bus 57 Random-No Limits" and "Bus57--With Limits")
# with IEEE Bus-57 shift-by-shift loss data (mirroring the existing
# GeoScenario sheet's layout), a new line chart on the first new sheet, and
# updates sheet/view selection state to match the authored workbook.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$geo = $wb.Worksheets.Item("GeoScenario")

# --- Add the two new worksheets, right after GeoScenario, in order ---------
$busRandom = $wb.Worksheets.Add($null, $geo)
$busRandom.Name = "Bus 57 Random-No Limits"

$busLimits = $wb.Worksheets.Add($null, $busRandom)
$busLimits.Name = "Bus57--With Limits"

# --- Populate "Bus 57 Random-No Limits" -------------------------------------
$busRandom.Range("A1").Value = "timestep"
$busRandom.Range("B1").Value = "RoadFirst"
$busRandom.Range("C1").Value = "Nominal Roads"
$busRandom.Range("D1").Value = "No Travel Times"

$dataRandom = @(
    @(0, 543.1, 543.1, 543.1),
    @(1, 159.29999999999899, 543.1, 159.29999999999899),
    @(2, 113.299999999999, 159.29999999999899, 113.299999999999),
    @(3, 91.3, 113.299999999999, 91.3),
    @(4, 88.6, 91.3, 91.3),
    @(5, 88.6, 88.3, 87.299999999999898)
)
for ($i = 0; $i -lt $dataRandom.Length; $i++) {
    $r = 2 + $i
    $row = $dataRandom[$i]
    $busRandom.Range("A$r").Value = $row[0]
    $busRandom.Range("B$r").Value = $row[1]
    $busRandom.Range("C$r").Value = $row[2]
    $busRandom.Range("D$r").Value = $row[3]
}

$busRandom.Range("A10").Value = "Summed"
for ($i = 0; $i -lt 6; $i++) {
    $r = 11 + $i
    $busRandom.Range("A$r").Value = $i
    $busRandom.Range("B$r").Formula = "=SUM(B2:B$(2+$i))"
    $busRandom.Range("C$r").Formula = "=SUM(C2:C$(2+$i))"
    $busRandom.Range("D$r").Formula = "=SUM(D2:D$(2+$i))"
}

$busRandom.Columns.Item(3).ColumnWidth = 14.42578125

# --- Populate "Bus57--With Limits" ------------------------------------------
$busLimits.Range("A1").Value = "timestep"
$busLimits.Range("B1").Value = "RoadFirst"
$busLimits.Range("C1").Value = "Nominal Roads"
$busLimits.Range("D1").Value = "No Travel Times"

$dataLimits = @(
    @(0, 517.1, 517.1, 517.1),
    @(1, 133.29999999999899, 133.29999999999899, 133.29999999999899),
    @(2, 87.299999999999898, 87.299999999999898, 87.299999999999898),
    @(3, 84.299999999999898, 83.699999999999903, 65.3),
    @(4, 83.699999999999903, 83.699999999999903, 52.3),
    @(5, 83.8, 83.6, 47.3)
)
for ($i = 0; $i -lt $dataLimits.Length; $i++) {
    $r = 2 + $i
    $row = $dataLimits[$i]
    $busLimits.Range("A$r").Value = $row[0]
    $busLimits.Range("B$r").Value = $row[1]
    $busLimits.Range("C$r").Value = $row[2]
    $busLimits.Range("D$r").Value = $row[3]
}

$busLimits.Range("A10").Value = "Summed"
for ($i = 0; $i -lt 6; $i++) {
    $r = 11 + $i
    $busLimits.Range("A$r").Value = $i
    $busLimits.Range("B$r").Formula = "=SUM(B2:B$(2+$i))"
    $busLimits.Range("C$r").Formula = "=SUM(C2:C$(2+$i))"
    $busLimits.Range("D$r").Formula = "=SUM(D2:D$(2+$i))"
}

$busLimits.Columns.Item(3).ColumnWidth = 14.42578125

# --- Chart on "Bus 57 Random-No Limits" -------------------------------------
$chartObjs = $busRandom.ChartObjects()
$chartObj = $chartObjs.Add(180975, 171450, 6200000, 3800000)
$chart = $chartObj.Chart
$chart.ChartType = 65
$chart.SetSourceData($busRandom.Range("A1:D7"))
$chart.HasTitle = $true
$chart.ChartTitle.Text = "Shift by Shift losses on IEEE Bus 57 without line limits"
$chart.HasLegend = $true
$chart.Legend.Position = -4107

# --- View / selection state --------------------------------------------------
$geo.Activate()
$geo.Range("A1:B16").Select()

$busRandom.Activate()
$busRandom.Range("A9:E17").Select()

$busLimits.Activate()
$busLimits.Range("H8").Select()
